$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (D) and 1h volume change (E) columns
# D-column cells whose new value is a plain number-like string are pre-formatted
# as Text so Excel keeps them as literal strings (matching the source inlineStr
# cells), then the style is reset to Normal so no stray formatting is left behind.

$ws.Range("D2").Value = '37.722.38'
$ws.Range("E2").Value = '  -0.19%  '
$ws.Range("D3").Value = '2.040.28'
$ws.Range("E3").Value = '  +0.35%  '
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '227.24'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.11%  '
$ws.Range("E6").Value = '  -0.91%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '59.47'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.94%  '
$ws.Range("E8").Value = '  -0.07%  '
$ws.Range("E9").Value = '  -2.65%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0838'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.52%  '
$ws.Range("E11").Value = '  +0.16%  '
$ws.Range("D12").Value = '2.340.84'
$ws.Range("E12").Value = '  +0.24%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.44'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.60%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '21.01'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.12%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.50'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +5.24%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.772'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.84%  '
$ws.Range("D17").Value = '2.043.18'
$ws.Range("E17").Value = '  +0.26%  '
$ws.Range("D18").Value = '37.712.14'
$ws.Range("E18").Value = '  -0.20%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '69.47'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.49%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.89'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.79%  '
$ws.Range("D21").Value = '0.0₃0822'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '223.55'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.90%  '
$ws.Range("E23").Value = '  +0.48%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.39'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.10%  '
$ws.Range("E25").Value = '  +2.32%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '168.61'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.16%  '
$ws.Range("E27").Value = '  +1.36%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.128'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.62%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '18.77'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.59%  '
$ws.Range("E30").Value = '  -0.03%  '
$ws.Range("E31").Value = '  -0.47%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.24'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +8.49%  '
$ws.Range("E33").Value = '  -1.13%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0603'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.22%  '
$ws.Range("E35").Value = '  +0.48%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.43'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.24%  '
$ws.Range("E37").Value = '  +3.70%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.42'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +5.07%  '
$ws.Range("E39").Value = '  +0.02%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '18.04'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +7.17%  '
$ws.Range("D41").Value = '1.526.33'
$ws.Range("E41").Value = '  -0.70%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '97.31'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.73%  '
$ws.Range("E43").Value = '  -1.00%  '
$ws.Range("E44").Value = '  +0.65%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.24'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +7.90%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0904'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.86%  '
$ws.Range("E47").Value = '  +0.03%  '
$ws.Range("E48").Value = '  +0.17%  '
$ws.Range("E49").Value = '  -0.91%  '
$ws.Range("E50").Value = '  -2.14%  '
$ws.Range("D51").Value = '2.230.67'
$ws.Range("E51").Value = '  +0.27%  '
